$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 219, shifting existing rows 219:325 down to 220:326
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new record
$ws.Cells.Item(219, 1).Value = 10
$ws.Cells.Item(219, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(219, 3).Value = "La Araucanía"
$ws.Cells.Item(219, 4).Value = 44523
$ws.Cells.Item(219, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(219, 5).Value = 9
$ws.Cells.Item(219, 6).Value = "Fruta"
$ws.Cells.Item(219, 7).Value = 100108
$ws.Cells.Item(219, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(219, 9).Value = 100108005
$ws.Cells.Item(219, 10).Value = "Piña"
$ws.Cells.Item(219, 11).Value = "Caramelo"
$ws.Cells.Item(219, 12).Value = "Primera"
$ws.Cells.Item(219, 13).Value = 85
$ws.Cells.Item(219, 14).Value = 19000
$ws.Cells.Item(219, 15).Value = 20000
$ws.Cells.Item(219, 16).Value = 19471
$ws.Cells.Item(219, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(219, 18).Value = "Ecuador"
$ws.Cells.Item(219, 19).Value = 1623
$ws.Cells.Item(219, 20).Value = 12
